# The presentation's single Design ("Integral" / "Red Violet") has its
# theme color scheme swapped for the plain default "Office Theme" colors
# (the scheme that, before this edit, only lived in the otherwise-unused
# theme part). PowerPoint's RGBColor.RGB setter on the master's
# ColorScheme is the supported COM surface for rewriting a theme's
# <a:clrScheme> color values, so drive the swap through that.

function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: the standard Office theme colors, in ColorScheme.Colors
# index order (1-12 == dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgbInt $officeThemeColors[$i - 1]
}
